$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells with same style as the existing header row (bold/bordered/centered)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the adjacent header cell (Unnamed: 28 / AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record (Wins/Losses/Ties) for every data row (2 through 47)
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 94
    $ws.Cells.Item($row, 31).Value = 68
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Output "done"
